$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "50XC" speed test result row (row 8):
# A8 already contains "Area, interp, + perim opt" and G8 already
# contains the commit hash; fill in the remaining columns.
$ws.Range("B8").Value = 2631
$ws.Range("C8").Value = 50
$ws.Range("D8").Value = "No"
$ws.Range("E8").Value = 10000
$ws.Range("F8").Value = "layered_multiXC_n50_fixed_ts_test.yml"

# Move the active selection, matching the recorded cursor position after edit.
$ws.Range("B13").Select()
